{"js": "// Apply the Harmonogram.docx edit:\n//  1) Swap the two dates between the \"-- wyszukiwanie lokali (po kryteriach )\"\n//     paragraph (25.11.2016 -> 2.12.2016) and the \"- komentowania\" paragraph\n//     (2.12.2016 -> 25.11.2016). The \"-ocenianie\" paragraph keeps 25.11.2016.\n//  2) Append, after the \"dodawanie edytowanie\" paragraph, a new blank\n//     paragraph followed by a paragraph describing the \"usuwanie\" (delete)\n//     confirmation screen work item.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet searchParagraph = null;\nlet koment = null;\nlet dodawanie = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"-- wyszukiwanie lokali\") !== -1) {\n    searchParagraph = paragraphs.items[i];\n  } else if (text.indexOf(\"- komentowania\") !== -1 || text.indexOf(\"\\t- komentowania\") !== -1) {\n    if (text.indexOf(\"System komentowania\") === -1) {\n      koment = paragraphs.items[i];\n    }\n  } else if (text.replace(/\\s+$/, \"\") === \"dodawanie edytowanie\") {\n    dodawanie = paragraphs.items[i];\n  }\n}\n\nif (!searchParagraph || !koment || !dodawanie) {\n  throw new Error(\n    \"Could not locate the expected paragraphs (\" +\n      \"search=\" + !!searchParagraph + \", koment=\" + !!koment + \", dodawanie=\" + !!dodawanie + \")\"\n  );\n}\n\n// 1a) \"-- wyszukiwanie lokali (po kryteriach )\" paragraph: 25.11.2016 -> 2.12.2016\nconst searchDate = searchParagraph.search(\"25.11.2016\", { matchCase: true });\nsearchDate.load(\"items\");\nawait context.sync();\nif (searchDate.items.length === 0) {\n  throw new Error(\"Could not find 25.11.2016 in the wyszukiwanie paragraph\");\n}\nsearchDate.items[0].insertText(\"2.12.2016\", \"Replace\");\nawait context.sync();\n\n// 1b) \"- komentowania\" paragraph: 2.12.2016 -> 25.11.2016\nconst komentDate = koment.search(\"2.12.2016\", { matchCase: true });\nkomentDate.load(\"items\");\nawait context.sync();\nif (komentDate.items.length === 0) {\n  throw new Error(\"Could not find 2.12.2016 in the komentowania paragraph\");\n}\nkomentDate.items[0].insertText(\"25.11.2016\", \"Replace\");\nawait context.sync();\n\n// 2) Append a blank paragraph and the new \"usuwanie\" paragraph at the end.\nconst blankPara = dodawanie.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nblankPara.insertParagraph(\n  \"usuwanie z customowym komunikatem czy na pewno chcemy usunac  zamienic przyciski tak nie miejscami na ca\u0142a strone, \u0142adne gui dla wyswietlania \",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "# Apply the Harmonogram.docx edit:\n#  1) Swap the two dates between the \"-- wyszukiwanie lokali (po kryteriach )\"\n#     paragraph (25.11.2016 -> 2.12.2016) and the \"- komentowania\" paragraph\n#     (2.12.2016 -> 25.11.2016). The \"-ocenianie\" paragraph keeps 25.11.2016.\n#  2) Append, after the \"dodawanie edytowanie\" paragraph, a new blank\n#     paragraph followed by a paragraph describing the \"usuwanie\" (delete)\n#     confirmation screen work item.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$searchPara = $null\n$komentPara = $null\n$dodawaniePara = $null\n\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"*wyszukiwanie lokali*\") {\n        $searchPara = $p\n    } elseif (($t -like \"*komentowania*\") -and ($t -notlike \"*System komentowania*\")) {\n        $komentPara = $p\n    } elseif ($t.Trim() -eq \"dodawanie edytowanie\") {\n        $dodawaniePara = $p\n    }\n}\n\nif (-not $searchPara) { throw \"Could not find the 'wyszukiwanie lokali' paragraph\" }\nif (-not $komentPara) { throw \"Could not find the 'komentowania' paragraph\" }\nif (-not $dodawaniePara) { throw \"Could not find the 'dodawanie edytowanie' paragraph\" }\n\n# 1a) \"-- wyszukiwanie lokali (po kryteriach )\" paragraph: 25.11.2016 -> 2.12.2016\n$find1 = $searchPara.Range.Find\n$find1.Text = \"25.11.2016\"\n$find1.Replacement.Text = \"2.12.2016\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 1b) \"- komentowania\" paragraph: 2.12.2016 -> 25.11.2016\n$find2 = $komentPara.Range.Find\n$find2.Text = \"2.12.2016\"\n$find2.Replacement.Text = \"25.11.2016\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 2) Append a blank paragraph and the new \"usuwanie\" paragraph at the end.\n$dodawaniePara.Range.InsertParagraphAfter()\n\n$count = $d.Paragraphs.Count\n$blankPara = $d.Paragraphs.Item($count)\n$blankPara.Range.InsertParagraphAfter()\n\n$count2 = $d.Paragraphs.Count\n$finalPara = $d.Paragraphs.Item($count2)\n$finalPara.Range.Text = \"usuwanie z customowym komunikatem czy na pewno chcemy usunac  zamienic przyciski tak nie miejscami na ca\u0142a strone, \u0142adne gui dla wyswietlania \"\n\nWrite-Output \"done\"\n"}
